$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '24.772.90'
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").Value = '1.682.03'
$ws.Range("E3").Value = '  -1.90%  '

Set-TextValue $ws.Range("D4") '1.006'
$ws.Range("E4").Value = '  +0.89%  '

Set-TextValue $ws.Range("D5") '313.61'
$ws.Range("E5").Value = '  -1.44%  '

Set-TextValue $ws.Range("D6") '1.005'
$ws.Range("E6").Value = '  +0.69%  '

Set-TextValue $ws.Range("D7") '0.3933'
$ws.Range("E7").Value = '  +0.03%  '

Set-TextValue $ws.Range("D8") '0.3966'
$ws.Range("E8").Value = '  -2.61%  '

$ws.Range("E9").Value = '  +1.08%  '

Set-TextValue $ws.Range("D10") '1.420'
$ws.Range("E10").Value = '  -5.40%  '

Set-TextValue $ws.Range("D11") '51.73'
$ws.Range("E11").Value = '  -3.61%  '

Set-TextValue $ws.Range("D12") '0.08668'
$ws.Range("E12").Value = '  -2.06%  '

$ws.Range("E13").Value = '  -4.57%  '

Set-TextValue $ws.Range("D14") '7.312'
$ws.Range("E14").Value = '  -2.95%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D15") '7.791'
$ws.Range("E15").Value = '  -4.35%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D16") '0.00001319'
$ws.Range("E16").Value = '  -3.41%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.569.96'
$ws.Range("E17").Value = '  -8.60%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D18") '94.08'
$ws.Range("E18").Value = '  -3.41%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D19") '0.07117'
$ws.Range("E19").Value = '  -1.42%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range("D20") '20.19'
$ws.Range("E20").Value = '  -2.83%  '

Set-TextValue $ws.Range("D21") '7.125'
$ws.Range("E21").Value = '  -2.87%  '

$ws.Range("E22").Value = '  +0.66%  '

Set-TextValue $ws.Range("D23") '14.09'
$ws.Range("E23").Value = '  -2.59%  '

$ws.Range("D24").Value = '24.779.81'
$ws.Range("E24").Value = '  -0.38%  '

Set-TextValue $ws.Range("D25") '2.352'
$ws.Range("E25").Value = '  +0.71%  '

Set-TextValue $ws.Range("D26") '23.74'
$ws.Range("E26").Value = '  +0.49%  '

Set-TextValue $ws.Range("D27") '2.773'
$ws.Range("E27").Value = '  -8.36%  '

Set-TextValue $ws.Range("D28") '162.35'
$ws.Range("E28").Value = '  -2.89%  '

$ws.Range("B29").Value = 'HuobiToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D29") '5.769'
$ws.Range("E29").Value = '  -3.46%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D30") '149.74'
$ws.Range("E30").Value = '  +2.17%  '

Set-TextValue $ws.Range("D31") '2.573'
$ws.Range("E31").Value = '  +14.19%  '

Set-TextValue $ws.Range("D32") '7.789'
$ws.Range("E32").Value = '  -9.10%  '

$ws.Range("B33").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C33").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D33").Value = '1.797.38'
$ws.Range("E33").Value = '  -5.69%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D34") '0.08431'
$ws.Range("E34").Value = '  -5.07%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D35") '0.03076'
$ws.Range("E35").Value = '  -2.89%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D36") '1.010'
$ws.Range("E36").Value = '  -4.86%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D37") '6.928'
$ws.Range("E37").Value = '  -5.08%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D38") '0.2802'
$ws.Range("E38").Value = '  -2.16%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D39") '0.09562'
$ws.Range("E39").Value = '  +3.16%  '

Set-TextValue $ws.Range("D40") '10.49'
$ws.Range("E40").Value = '  -4.31%  '

Set-TextValue $ws.Range("D41") '0.7943'
$ws.Range("E41").Value = '  -6.45%  '

$ws.Range("E42").Value = '  -0.91%  '

Set-TextValue $ws.Range("D43") '13.68'
$ws.Range("E43").Value = '  -3.66%  '

Set-TextValue $ws.Range("D44") '16.63'
$ws.Range("E44").Value = '  -4.97%  '

Set-TextValue $ws.Range("D45") '0.7164'
$ws.Range("E45").Value = '  -4.40%  '

Set-TextValue $ws.Range("D46") '2.580'
$ws.Range("E46").Value = '  -4.87%  '

Set-TextValue $ws.Range("D47") '4.181'
$ws.Range("E47").Value = '  -2.24%  '

$ws.Range("E48").Value = '  +4.87%  '

Set-TextValue $ws.Range("D49") '1.004'
$ws.Range("E49").Value = '  +0.58%  '

Set-TextValue $ws.Range("D50") '1.340'
$ws.Range("E50").Value = '  -4.93%  '

Set-TextValue $ws.Range("D51") '138.14'
$ws.Range("E51").Value = '  -1.97%  '
